{"js": "// 1) \"shocked\" -> \"tickled\" in the sentence about Mara's tsundere behaviour.\nconst shockedResults = context.document.body.search(\"a little shocked by Mara\\u2019s sudden tsundere behaviour\", { matchCase: true });\nshockedResults.load(\"text\");\nawait context.sync();\n\nif (shockedResults.items.length > 0) {\n  shockedResults.items[0].insertText(\n    \"a little tickled by Mara\\u2019s sudden tsundere behaviour\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 2) Consolidate the three runs that make up the \"dozes off, contently clutching\n//    the jacket I left on the couch.\" sentence into a single run. The visible\n//    text is unchanged; re-writing the whole sentence as one insertText call\n//    causes the underlying runs to merge into one.\nconst fullSentence =\n  \"We end up watching the accompanying movie to one of the shows we loved as kids, but halfway through Mara dozes off, contently clutching the jacket I left on the couch.\";\n\nconst sentenceResults = context.document.body.search(fullSentence, { matchCase: true });\nsentenceResults.load(\"text\");\nawait context.sync();\n\nif (sentenceResults.items.length > 0) {\n  sentenceResults.items[0].insertText(fullSentence, \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"shocked\" -> \"tickled\" in the sentence about Mara's tsundere behaviour.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"shocked\"\n$find1.Replacement.Text = \"tickled\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) Consolidate the three runs that make up the \"dozes off, contently clutching\n#    the jacket I left on the couch.\" sentence into a single run. The visible\n#    text does not change; re-running Find & Replace across the whole sentence\n#    (even though the replacement text equals the existing text) causes Word to\n#    rewrite the matched range as one run instead of the original three.\n$fullText = \"We end up watching the accompanying movie to one of the shows we loved as kids, but halfway through Mara dozes off, contently clutching the jacket I left on the couch.\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = $fullText\n$find2.Replacement.Text = $fullText\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
